$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D ("Price") and E ("Volume(1h)") hold plain text values
# (e.g. "66.076.71", "  -0.79%  "), not numbers. Prefix each new value
# with a leading apostrophe (Excel's text-literal marker) so numeric-
# looking strings such as "585.82" are written back as text instead of
# being reinterpreted as numbers (matches the original cell typing).

$ws.Range("D2").Value = "'66.076.71"
$ws.Range("E2").Value = "'  -0.79%  "

$ws.Range("D3").Value = "'3.315.42"
$ws.Range("E3").Value = "'  -0.57%  "

$ws.Range("E4").Value = "'  +0.00%  "

$ws.Range("D5").Value = "'585.82"
$ws.Range("E5").Value = "'  +2.15%  "

$ws.Range("D6").Value = "'181.93"
$ws.Range("E6").Value = "'  +0.86%  "

$ws.Range("D7").Value = "'0.653"
$ws.Range("E7").Value = "'  +3.77%  "

$ws.Range("E8").Value = "'  +0.02%  "

$ws.Range("D9").Value = "'3.314.03"
$ws.Range("E9").Value = "'  -0.59%  "

$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = "'  -3.03%  "

$ws.Range("D11").Value = "'6.81"
$ws.Range("E11").Value = "'  +2.56%  "

$ws.Range("D12").Value = "'0.401"
$ws.Range("E12").Value = "'  -0.57%  "

$ws.Range("D13").Value = "'3.892.02"
$ws.Range("E13").Value = "'  -0.59%  "

$ws.Range("E14").Value = "'  -2.93%  "

$ws.Range("D15").Value = "'66.148.50"
$ws.Range("E15").Value = "'  -0.87%  "

$ws.Range("D16").Value = "'26.13"
$ws.Range("E16").Value = "'  -3.18%  "

$ws.Range("E17").Value = "'  -1.22%  "

$ws.Range("D18").Value = "'3.278.08"
$ws.Range("E18").Value = "'  -2.21%  "

$ws.Range("D19").Value = "'424.71"
$ws.Range("E19").Value = "'  -2.77%  "

$ws.Range("D20").Value = "'5.54"
$ws.Range("E20").Value = "'  -2.52%  "

$ws.Range("D21").Value = "'13.11"
$ws.Range("E21").Value = "'  -3.09%  "

$ws.Range("D22").Value = "'7.37"
$ws.Range("E22").Value = "'  -2.85%  "

$ws.Range("D23").Value = "'71.68"
$ws.Range("E23").Value = "'  -2.49%  "

$ws.Range("E24").Value = "'  -0.07%  "

$ws.Range("E25").Value = "'  +0.33%  "

$ws.Range("D26").Value = "'3.460.75"
$ws.Range("E26").Value = "'  -0.65%  "

$ws.Range("D27").Value = "'0.512"
$ws.Range("E27").Value = "'  -0.87%  "

$ws.Range("D28").Value = "'0.200"
$ws.Range("E28").Value = "'  +4.91%  "

$ws.Range("D29").Value = "'0.0000113"
$ws.Range("E29").Value = "'  -3.13%  "

$ws.Range("D30").Value = "'8.87"
$ws.Range("E30").Value = "'  -1.74%  "

$ws.Range("E31").Value = "'  +0.21%  "

$ws.Range("E32").Value = "'  -2.14%  "

$ws.Range("D33").Value = "'22.36"
$ws.Range("E33").Value = "'  -1.90%  "

$ws.Range("E35").Value = "'  -1.82%  "

$ws.Range("D36").Value = "'6.55"
$ws.Range("E36").Value = "'  -3.07%  "

$ws.Range("E37").Value = "'  -4.68%  "

$ws.Range("D38").Value = "'160.58"
$ws.Range("E38").Value = "'  -1.75%  "

$ws.Range("E39").Value = "'  -3.31%  "

$ws.Range("D40").Value = "'2.870.75"
$ws.Range("E40").Value = "'  +1.45%  "

$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "'  +0.20%  "

$ws.Range("D42").Value = "'26.34"
$ws.Range("E42").Value = "'  -3.54%  "

$ws.Range("E43").Value = "'  -4.70%  "

$ws.Range("D44").Value = "'4.30"
$ws.Range("E44").Value = "'  -2.61%  "

$ws.Range("D45").Value = "'39.81"
$ws.Range("E45").Value = "'  -1.05%  "

$ws.Range("D46").Value = "'0.0659"
$ws.Range("E46").Value = "'  -1.16%  "

$ws.Range("D47").Value = "'5.90"
$ws.Range("E47").Value = "'  -4.94%  "

$ws.Range("D48").Value = "'2.28"
$ws.Range("E48").Value = "'  -2.65%  "

$ws.Range("D49").Value = "'23.08"
$ws.Range("E49").Value = "'  -5.36%  "

$ws.Range("D50").Value = "'312.57"
$ws.Range("E50").Value = "'  -2.91%  "

$ws.Range("E51").Value = "'  -0.95%  "

